# Apply the "Added more county data" commit to the Adjacency_matrix sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Adjacency_matrix")

# 1. Rename the sheet (Adjacency_matrix -> Coiunty)
$ws.Name = "Coiunty"

# 2. Add the new "County" label in A1 and "Pop_2016" header in AB1
$ws.Cells.Item(1, 1).Value = "County"
$ws.Cells.Item(1, 28).Value = "Pop_2016"

# 3. Populate the new Pop_2016 column (AB) for each county row (2-27)
#    Values are the 2016 Census populations, in the same row order as the
#    existing adjacency matrix.
$pop2016 = @{
    2  = 56932    # Carlow
    3  = 76176    # Cavan
    4  = 118817   # Clare
    5  = 542868   # Cork
    6  = 159192   # Donegal
    7  = 1347359  # Dublin
    8  = 258058   # Galway
    9  = 147707   # Kerry
    10 = 222504   # Kildare
    11 = 99232    # Kilkenny
    12 = 84697    # Laois
    13 = 32044    # Leitrim
    14 = 194899   # Limerick
    15 = 40873    # Longford
    16 = 128884   # Louth
    17 = 130507   # Mayo
    18 = 195044   # Meath
    19 = 61386    # Monaghan
    20 = 77961    # Offaly
    21 = 64544    # Roscommon
    22 = 65535    # Sligo
    23 = 159553   # Tipperary
    24 = 116176   # Waterford
    25 = 88770    # Westmeath
    26 = 149722   # Wexford
    27 = 142425   # Wicklow
}

foreach ($r in 2..27) {
    $ws.Cells.Item($r, 28).Value = $pop2016[$r]
}

# 4. Append a block of mostly-empty, number-formatted cells in columns
#    U (#,##0), V (#,##0.00, only on a few rows) and X (0.00%) for rows 31-65.
$vRows = @(39, 40, 59)

foreach ($r in 31..65) {
    $ws.Cells.Item($r, 21).NumberFormat = "#,##0"      # column U
    if ($vRows -contains $r) {
        $ws.Cells.Item($r, 22).NumberFormat = "#,##0.00"  # column V
    }
    if ($r -ne 37 -and $r -ne 53 -and $r -ne 58) {
        $ws.Cells.Item($r, 24).NumberFormat = "0.00%"      # column X
    }
}

# 5. Restore the view state (scrolled pane / selected cell) recorded for
#    this sheet after the edits were made.
$ws.Range("G11").Select()

Write-Host "done"
